# Update "want-to-go count" (F column) figures on the 展览 and 全部类型 sheets
# to reflect freshly generated output (gh-pages build 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (rows keyed by their row number -> new F value)
$exhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    2  = 1316
    3  = 1187
    4  = 14440
    5  = 16979
    7  = 132
    8  = 44
    16 = 40
    19 = 1295
    22 = 56
    23 = 37
    24 = 5
    25 = 6985
    27 = 29
    28 = 1149
    29 = 18
    31 = 5809
    34 = 208
    35 = 4952
}
foreach ($row in $exhibitionUpdates.Keys) {
    $exhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (row numbers shifted by extra rows present in this sheet)
$allTypes = $wb.Worksheets.Item("全部类型")
$allTypesUpdates = @{
    2  = 1316
    3  = 1187
    4  = 14440
    5  = 16979
    7  = 132
    8  = 44
    16 = 40
    19 = 1295
    23 = 56
    24 = 37
    25 = 5
    26 = 6985
    28 = 29
    29 = 1149
    30 = 18
    33 = 5809
    36 = 208
    37 = 4952
}
foreach ($row in $allTypesUpdates.Keys) {
    $allTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
